$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that sits in the middle of the
#    abstract paragraph (between " their" and " IoT network which can help...").
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# 2. Merge the two survey-question paragraphs
#       "have" + " you upgrade your IoT devices?"
#       "Have you change default setting on IoT devices?"
#    into a single corrected paragraph:
#       "Have you changed the default setting on IoT devices?"
#    Locate the "have you upgrade ..." paragraph and delete it (together with
#    its paragraph mark), which merges what follows it into the next
#    paragraph; then fix up the wording of the remaining paragraph.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "have you upgrade*") {
        $targetIndex = $i
    }
    $i++
}
if ($targetIndex -ge 1) {
    $d.Paragraphs($targetIndex).Range.Delete()
}

$null = $d.Content.Find.Execute("Have you change default setting on IoT devices?", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Have you changed the default setting on IoT devices?", 2)

# 3. Append two more empty paragraphs at the very end of the document
#    (there were already two trailing empty paragraphs). The last of the
#    two new paragraphs holds a fresh "_GoBack" bookmark, matching where
#    Word itself drops that marker after the most recent edit/save.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
